$d = $word.ActiveDocument

# 1. Remove the entire "Add parentheses logic" bullet paragraph.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Add paren*theses logic*") {
        $p.Range.Delete()
        break
    }
}

# 2. Trim the trailing "if I have the knowledge" from the final bullet.
$d.Content.Find.Execute(
    "Review all incomplete items listed here and complete if I have the knowledge",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Review all incomplete items listed here and complete",
    2)
